{"js": "// Title paragraph edit for \"Qualified Event Names\":\n//   1. Remove the (hidden) \"_GoBack\" bookmark that wraps the title run.\n//   2. Change the title text from\n//        \"Circle Language Spec: Events\"\n//      to\n//        \"Circle Language Construct Drafts | Events\"\n\nconst OLD_TITLE = \"Circle Language Spec: Events\";\nconst NEW_TITLE = \"Circle Language Construct Drafts | Events\";\n\n// --- Step 1: drop the \"_GoBack\" bookmark wrapping the title paragraph ---\n// The Word.js bookmark-deletion helpers (Document.deleteBookmark /\n// Range.delete on a bookmark range) are not reliable for this hidden\n// bookmark in this host, so rebuild the paragraph's OOXML without the\n// <w:bookmarkStart>/<w:bookmarkEnd> pair and replace the paragraph with it.\n// This preserves every other attribute (rsids, pStyle, run contents, \u2026).\nconst titleResults = context.document.body.search(OLD_TITLE, { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  const titlePara = titleResults.items[0].paragraphs.getFirst();\n  const titleRange = titlePara.getRange();\n  const ooxml = titlePara.getOoxml();\n  await context.sync();\n\n  let pXml = ooxml.value.match(/<w:p(?:\\s[^>]*)?>[\\s\\S]*?<\\/w:p>/)[0];\n\n  // Strip the bookmark markers.\n  pXml = pXml.replace(/<w:bookmarkStart[^>]*\\/>/g, \"\");\n  pXml = pXml.replace(/<w:bookmarkEnd[^>]*\\/>/g, \"\");\n  // getOoxml() mints fresh w14:paraId/w14:textId values for round-tripping;\n  // the source paragraph never had them, so drop them again.\n  pXml = pXml.replace(/\\s+w14:paraId=\"[^\"]*\"/g, \"\");\n  pXml = pXml.replace(/\\s+w14:textId=\"[^\"]*\"/g, \"\");\n\n  const packageXml =\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    pXml +\n    \"</w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n  titleRange.insertOoxml(packageXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Step 2: update the title text itself ---\nconst searchResults = context.document.body.search(OLD_TITLE, { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(NEW_TITLE, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Title paragraph edit for \"Qualified Event Names\":\n#   1. Remove the (hidden) \"_GoBack\" bookmark that wraps the title run.\n#   2. Change the title text from\n#        \"Circle Language Spec: Events\"\n#      to\n#        \"Circle Language Construct Drafts | Events\"\n\n$d = $word.ActiveDocument\n\n# --- Step 1: drop the \"_GoBack\" bookmark wrapping the title paragraph ---\n# The Bookmarks collection hides \"_GoBack\" (a hidden/system bookmark) from\n# enumeration/Delete in this host, so rebuild the first paragraph's OOXML\n# without the <w:bookmarkStart>/<w:bookmarkEnd> pair and re-insert it in\n# place. This preserves every other attribute (rsids, pStyle, run contents, \u2026).\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$fullXml = $titleRange.XML()\n\n$match = [regex]::Match($fullXml, '<w:p(?:\\s[^>]*)?>[\\s\\S]*?</w:p>')\nif ($match.Success) {\n    $pXml = $match.Value\n\n    # Strip the bookmark markers.\n    $pXml = [regex]::Replace($pXml, '<w:bookmarkStart[^>]*/>', '')\n    $pXml = [regex]::Replace($pXml, '<w:bookmarkEnd[^>]*/>', '')\n    # .XML() mints fresh w14:paraId/w14:textId values for round-tripping;\n    # the source paragraph never had them, so drop them again.\n    $pXml = [regex]::Replace($pXml, '\\s+w14:paraId=\"[^\"]*\"', '')\n    $pXml = [regex]::Replace($pXml, '\\s+w14:textId=\"[^\"]*\"', '')\n\n    $titleRange.InsertXML($pXml)\n}\n\n# --- Step 2: update the title text itself ---\n$find = $d.Content.Find\n$find.Text = \"Circle Language Spec: Events\"\n$find.Replacement.Text = \"Circle Language Construct Drafts | Events\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
